$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A8 with new consolidated values
$ws.Range("A2").Value = "('Bear', ['Token Creature — Bear', '2/2'])"
$ws.Range("A3").Value = "('Demon', ['Token Creature — Demon', 'Flying', '*/*'])"
$ws.Range("A4").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A5").Value = "('Insect', ['Token Creature — Insect', '1/1'])"
$ws.Range("A6").Value = "('Rukh', ['Token Creature — Rukh', 'Flying', '4/4'])"
$ws.Range("A7").Value = "('Sliver', ['Token Creature — Sliver', '1/1'])"
$ws.Range("A8").Value = "('Voidmage Prodigy', ['{U}{U}', 'Creature — Human Wizard', '{U}{U}, Sacrifice a Wizard: Counter target spell.', 'Morph {U} (You may cast this card face down as a 2/2 creature for {3}. Turn it face up any time for its morph cost.)', '2/1'])"

# Delete the now-obsolete rows 9 through 27
$ws.Range("A9:A27").EntireRow.Delete()
